$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Translate header "Año" -> "YEAR" (rate stays the same)
$ws.Range("A1").Value = "YEAR"
$ws.Range("B1").Value = "rate"

# Move the active selection (cosmetic, matches recorded cursor position)
$ws.Range("E7").Select()
